# Apply updated odds values to Sheet1 (cells previously holding numeric odds data)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F2").Value = 1.65
$ws.Range("G2").Value = 1.78
$ws.Range("H2").Value = 5.7
$ws.Range("I2").Value = 6.8
$ws.Range("J2").Value = 3.75
$ws.Range("K2").Value = 4.4
$ws.Range("M2").Value = 1.07
$ws.Range("N2").Value = 3.45
$ws.Range("O2").Value = 1.32
$ws.Range("Q2").Value = 1.88
$ws.Range("U2").Value = 1.9
$ws.Range("V2").Value = 1.16
$ws.Range("W2").Value = 2.38
$ws.Range("AB2").Value = 9.800000000000001
$ws.Range("AC2").Value = 11.5
$ws.Range("AF2").Value = 11.5
$ws.Range("AG2").Value = 12
$ws.Range("H3").Value = 2.34
$ws.Range("X3").Value = 19
$ws.Range("Z3").Value = 21
$ws.Range("AA3").Value = 42
$ws.Range("AF3").Value = 26
$ws.Range("AK3").Value = 42
$ws.Range("AN3").Value = 36
$ws.Range("G4").Value = 1.44
$ws.Range("L4").Value = 1.01
$ws.Range("R4").Value = 1.5
$ws.Range("S4").Value = 2.22
$ws.Range("G5").Value = 3.7
$ws.Range("H5").Value = 2.28
$ws.Range("J5").Value = 3.4
$ws.Range("K5").Value = 4
$ws.Range("M5").Value = 1.05
$ws.Range("N5").Value = 3.7
$ws.Range("O5").Value = 1.2
$ws.Range("P5").Value = 2
$ws.Range("Q5").Value = 1.79
$ws.Range("R5").Value = 1.11
$ws.Range("S5").Value = 2.66
$ws.Range("W5").Value = 1.37
$ws.Range("Z5").Value = 23
$ws.Range("AA5").Value = 46
$ws.Range("AD5").Value = 16
$ws.Range("AF5").Value = 32
$ws.Range("AG5").Value = 19
$ws.Range("AH5").Value = 23
$ws.Range("H6").Value = 3.4
$ws.Range("I6").Value = 3.95
$ws.Range("L6").Value = 1.01
$ws.Range("M6").Value = 1.01
$ws.Range("N6").Value = 1.63
$ws.Range("O6").Value = 1.44
$ws.Range("P6").Value = 1.63
$ws.Range("R6").Value = 1.08
$ws.Range("S6").Value = 3.75
$ws.Range("T6").Value = 1.01
$ws.Range("U6").Value = 1.01
$ws.Range("V6").Value = 1.37
$ws.Range("W6").Value = 1.64
$ws.Range("X6").Value = 14.5
$ws.Range("Y6").Value = 15.5
$ws.Range("Z6").Value = 34
$ws.Range("AA6").Value = 100
$ws.Range("AB6").Value = 11
$ws.Range("AC6").Value = 10
$ws.Range("AD6").Value = 21
$ws.Range("AE6").Value = 70
$ws.Range("AF6").Value = 20
$ws.Range("AG6").Value = 16
$ws.Range("AH6").Value = 29
$ws.Range("AI6").Value = 95
$ws.Range("AJ6").Value = 48
$ws.Range("AK6").Value = 44
$ws.Range("AL6").Value = 70
$ws.Range("AM6").Value = 1000
$ws.Range("AN6").Value = 1000
$ws.Range("AO6").Value = 1000
$ws.Range("F7").Value = 1.43
$ws.Range("I7").Value = 12.5
$ws.Range("O7").Value = 1.45
$ws.Range("Q7").Value = 2.34
$ws.Range("V7").Value = 1.08
$ws.Range("AG7").Value = 13.5
$ws.Range("AI7").Value = 300
$ws.Range("G11").Value = 2.48
$ws.Range("J11").Value = 2.98
$ws.Range("Q11").Value = 2.44
$ws.Range("P13").Value = 1.99
